$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.521.43"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "3.825.63"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'701.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.25%  "
$ws.Range("D6").Value = "'174.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("D7").Value = "3.824.73"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("E11").Value = "  +5.88%  "
$ws.Range("D12").Value = "'0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("E13").Value = "  +6.25%  "
$ws.Range("D14").Value = "'36.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "4.462.23"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "3.811.26"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "71.497.56"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "'17.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "'7.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D21").Value = "'11.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("D22").Value = "'486.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.42%  "
$ws.Range("D23").Value = "'0.718"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").Value = "'84.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("D25").Value = "'0.0000144"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").Value = "'10.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("D29").Value = "3.972.78"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "'3.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.03%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  +2.73%  "
$ws.Range("D34").Value = "'29.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("D36").Value = "'9.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.777.10"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("B38").Value = "Binance-PegBSC-USD"
$ws.Range("C38").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D38").Value = "'0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.104"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.30%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'6.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'0.000311"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.57%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'163.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "'44.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'48.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").Value = "'0.303"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("E51").Value = "  +2.47%  "
